$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SDRF")

# The combined "ParameterValue" column (Z) is being split back out: the
# original single column is removed (its former neighbour slides into
# its place) and every later column shifts one place to the left.
$ws.Columns.Item(26).Delete()

# After the shift, the columns that used to be AW/AX are now AV/AW.
# AV keeps the "RUM pipeline: alignment and coverage" text but each row
# gets its own "<id=N>" suffix; AW keeps its text ("Elisabetta Manduchi")
# but picks up the yellow highlight fill that used to live on the
# deleted column (while keeping its own font).
$ids = 1,2,3,4,5,6
for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $ws.Range("AV$row").Value = "RUM pipeline: alignment and coverage <id=$($ids[$i])>"
    $ws.Range("AW$row").Interior.Color = 65535
}
